# Commit: "remove unused define for item"
#
# The ItemSubType column (C) held leftover per-row "define" values that are
# no longer used; clear them to 0. For the block of rows that represented a
# single item definition repeated across sub-levels (rows 46-61) the Level
# column (B) is also renumbered from 5 to 2 while ItemSubType becomes 7.
# Finally, restore the sheet's persisted view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Rows 13-45: ItemSubType -> 0 (Level in column B is untouched).
$ws.Range("C13:C45").Value = 0

# Rows 46-61: Level 5 -> 2, ItemSubType 0 -> 7.
$ws.Range("B46:B61").Value = 2
$ws.Range("C46:C61").Value = 7

# Rows 132-140: ItemSubType -> 0.
$ws.Range("C132:C140").Value = 0

# Restore the saved view/selection state (frozen-pane scroll + active cell).
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 125
$window.ScrollColumn = 2
$ws.Range("D131").Select()
